$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 70, shifting existing rows 70-118 down to 71-119
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new record
$ws.Cells.Item(70, 1).Value = 4
$ws.Cells.Item(70, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(70, 3).Value = "Los Lagos"
$ws.Cells.Item(70, 4).Value = 44729
$ws.Cells.Item(70, 5).Value = 10
$ws.Cells.Item(70, 6).Value = 100112022
$ws.Cells.Item(70, 7).Value = "Arveja Verde"
$ws.Cells.Item(70, 8).Value = "Perfection"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 35
$ws.Cells.Item(70, 11).Value = 44000
$ws.Cells.Item(70, 12).Value = 44000
$ws.Cells.Item(70, 13).Value = 44000
$ws.Cells.Item(70, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(70, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(70, 16).Value = 1760
$ws.Cells.Item(70, 17).Value = 25
$ws.Cells.Item(70, 18).Value = "Hortaliza"

# Match the date style used by column D in the surrounding rows
$ws.Cells.Item(70, 4).NumberFormat = $ws.Cells.Item(71, 4).NumberFormat
